# SA-16 Israel Time SA-17 Forbes on Russian SA-18 Coindesk The Block Crypto
#
# Categories.xlsx has two sheets:
#   COMMON  (categories)            -> columns: A=name, B=parent, C=ru, D=en, E=ukr
#   REGION  (regions / countries)   -> columns: A=name, B=parent, C=ru, D=en, E=ukr
#
# This change:
#   1) Fixes mismatched translations for the Theatre / Science / History /
#      Animals rows on COMMON (rows 51-54), which had been pointing at each
#      other's ru/en/ukr strings.
#   2) Adds a new "comedy" category (child of Culture) on COMMON.
#   3) Adds a new "israel" region (child of MiddleEast) on REGION.

$wb = $excel.ActiveWorkbook

$common = $wb.Worksheets.Item("COMMON")
$region = $wb.Worksheets.Item("REGION")

# --- 1) Fix rows 51-54 on COMMON (Theatre / Science / History / Animals) ---

# Row 51: Theatre
$common.Range("A51").Value = "Theatre"
$common.Range("B51").Value = "Culture"
$common.Range("C51").Value = "Театр"
$common.Range("D51").Value = "Theatre"
$common.Range("E51").Value = "Театр"

# Row 52: Science
$common.Range("A52").Value = "Science"
$common.Range("C52").Value = "Наука"
$common.Range("D52").Value = "Science"
$common.Range("E52").Value = "Наука"

# Row 53: History
$common.Range("A53").Value = "History"
$common.Range("B53").Value = "Culture"
$common.Range("C53").Value = "История"
$common.Range("D53").Value = "History"
$common.Range("E53").Value = "Історія"

# Row 54: Animals
$common.Range("A54").Value = "Animals"
$common.Range("B54").Value = "Ecology"
$common.Range("C54").Value = "Животные"
$common.Range("D54").Value = "Animals"
$common.Range("E54").Value = "Тварини"

# --- 2) Add new "comedy" category (row 85) on COMMON ---

$common.Range("A85").Value = "comedy"
$common.Range("B85").Value = "Culture"
$common.Range("C85").Value = "Комедия"
$common.Range("D85").Value = "Comedy"
$common.Range("E85").Value = "Комедія"

# --- 3) Add new "israel" region (row 54) on REGION ---

$region.Range("A54").Value = "israel"
$region.Range("B54").Value = "MiddleEast"
$region.Range("C54").Value = "Израиль"
$region.Range("D54").Value = "Israel"
$region.Range("E54").Value = "Ізраїль"
